# Auto-generated Excel COM-interop edit script.
# Applies numeric corrections to columns H-N across multiple sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the scheduled-runner update.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 8575.362999999999
$ws.Range("I6").Value = 8707.571
$ws.Range("J6").Value = 5799
$ws.Range("K6").Value = 26122.713
$ws.Range("L6").Value = 17397
$ws.Range("M6").Value = -26010.713
$ws.Range("N6").Value = -17621
$ws.Range("H18").Value = 349
$ws.Range("I18").Value = 349
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 349
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -65
$ws.Range("N18").Value = $null
$ws.Range("H28").Value = 706.3077
$ws.Range("I28").Value = 389.27274
$ws.Range("K28").Value = 389.27274
$ws.Range("M28").Value = 95.72726
$ws.Range("H34").Value = 11136
$ws.Range("I34").Value = 11136
$ws.Range("K34").Value = 11136
$ws.Range("M34").Value = -10933
$ws.Range("H36").Value = 11136
$ws.Range("I36").Value = 11136
$ws.Range("K36").Value = 11136
$ws.Range("M36").Value = -10421
$ws.Range("H106").Value = 18365.709
$ws.Range("I106").Value = 5386.2
$ws.Range("K106").Value = 5386.2
$ws.Range("M106").Value = -4755.2
$ws.Range("H107").Value = 1050.0435
$ws.Range("I107").Value = 861.7646999999999
$ws.Range("K107").Value = 861.7646999999999
$ws.Range("M107").Value = 1058.2353
$ws.Range("H125").Value = 1050
$ws.Range("I125").Value = 1250
$ws.Range("K125").Value = 11250
$ws.Range("M125").Value = -8790
$ws.Range("H137").Value = 11361.728
$ws.Range("I137").Value = 3897.4827
$ws.Range("J137").Value = 25792.6
$ws.Range("K137").Value = 11692.4481
$ws.Range("L137").Value = 77377.79999999999
$ws.Range("M137").Value = -9142.4481
$ws.Range("N137").Value = -82477.79999999999

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4413.1763
$ws.Range("J2").Value = 8137.2856
$ws.Range("L2").Value = 8137.2856
$ws.Range("N2").Value = -8363.285599999999
$ws.Range("H97").Value = 738.4286
$ws.Range("I97").Value = 695.30304
$ws.Range("K97").Value = 695.30304
$ws.Range("M97").Value = -199.30304
$ws.Range("H102").Value = 1899.9524
$ws.Range("I102").Value = 1627.7222
$ws.Range("K102").Value = 1627.7222
$ws.Range("M102").Value = -5.72219999999993
$ws.Range("H116").Value = 4413.1763
$ws.Range("J116").Value = 8137.2856
$ws.Range("L116").Value = 8137.2856
$ws.Range("N116").Value = -12725.2856
$ws.Range("H122").Value = 2935.4285
$ws.Range("I122").Value = 3078.6
$ws.Range("J122").Value = 2577.5
$ws.Range("K122").Value = 9235.799999999999
$ws.Range("L122").Value = 7732.5
$ws.Range("M122").Value = -6785.799999999999
$ws.Range("N122").Value = -12632.5

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4413.1763
$ws.Range("J3").Value = 8137.2856
$ws.Range("L3").Value = 8137.2856
$ws.Range("N3").Value = -8365.285599999999
$ws.Range("H94").Value = 336.70587
$ws.Range("I94").Value = 248.93333
$ws.Range("J94").Value = 995
$ws.Range("K94").Value = 248.93333
$ws.Range("L94").Value = 995
$ws.Range("M94").Value = 202.06667
$ws.Range("N94").Value = -1897
$ws.Range("H105").Value = 3047.4
$ws.Range("I105").Value = 3342.4707
$ws.Range("K105").Value = 3342.4707
$ws.Range("M105").Value = -1595.4707
$ws.Range("H107").Value = 2694.875
$ws.Range("I107").Value = 2777.4614
$ws.Range("J107").Value = 2337
$ws.Range("K107").Value = 2777.4614
$ws.Range("L107").Value = 2337
$ws.Range("M107").Value = -857.4614000000001
$ws.Range("N107").Value = -6177

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 19463
$ws.Range("J41").Value = 21644.584
$ws.Range("L41").Value = 21644.584
$ws.Range("N41").Value = -22500.584
$ws.Range("H51").Value = 18398.334
$ws.Range("I51").Value = 18398.334
$ws.Range("K51").Value = 18398.334
$ws.Range("M51").Value = -17662.334
$ws.Range("H60").Value = 9046.333000000001
$ws.Range("I60").Value = 9046.333000000001
$ws.Range("K60").Value = 9046.333000000001
$ws.Range("M60").Value = -8535.333000000001
$ws.Range("H61").Value = 18398.334
$ws.Range("I61").Value = 18398.334
$ws.Range("K61").Value = 18398.334
$ws.Range("M61").Value = -18050.334

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 406.42856
$ws.Range("I92").Value = 369.8
$ws.Range("K92").Value = 1109.4
$ws.Range("M92").Value = 138.5999999999999
$ws.Range("H132").Value = 33340000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 33340000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 300060000
$ws.Range("M132").Value = $null
$ws.Range("N132").Value = -300065060

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 26922.666
$ws.Range("I58").Value = 12346.667
$ws.Range("J58").Value = 34210.668
$ws.Range("K58").Value = 12346.667
$ws.Range("L58").Value = 34210.668
$ws.Range("M58").Value = -12069.667
$ws.Range("N58").Value = -34764.668
$ws.Range("H113").Value = 156560.53
$ws.Range("I113").Value = 184490.19
$ws.Range("K113").Value = 184490.19
$ws.Range("M113").Value = -182320.19
$ws.Range("H126").Value = 3036.9092
$ws.Range("I126").Value = 2372.8572
$ws.Range("J126").Value = 4199
$ws.Range("K126").Value = 7118.571599999999
$ws.Range("L126").Value = 12597
$ws.Range("M126").Value = -4648.571599999999
$ws.Range("N126").Value = -17537

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2802.2856
$ws.Range("I16").Value = 1602.6666
$ws.Range("J16").Value = 10000
$ws.Range("K16").Value = 1602.6666
$ws.Range("L16").Value = 10000
$ws.Range("M16").Value = -1432.6666
$ws.Range("N16").Value = -10340
$ws.Range("H22").Value = 2464.64
$ws.Range("I22").Value = 1807.5
$ws.Range("J22").Value = 3071.2307
$ws.Range("K22").Value = 1807.5
$ws.Range("L22").Value = 3071.2307
$ws.Range("M22").Value = -1512.5
$ws.Range("N22").Value = -3661.2307
$ws.Range("H27").Value = 2464.64
$ws.Range("I27").Value = 1807.5
$ws.Range("J27").Value = 3071.2307
$ws.Range("K27").Value = 1807.5
$ws.Range("L27").Value = 3071.2307
$ws.Range("M27").Value = -1700.5
$ws.Range("N27").Value = -3285.2307
$ws.Range("H46").Value = 1025.9445
$ws.Range("I46").Value = 1061.375
$ws.Range("J46").Value = 997.6
$ws.Range("K46").Value = 1061.375
$ws.Range("L46").Value = 997.6
$ws.Range("M46").Value = -873.375
$ws.Range("N46").Value = -1373.6
$ws.Range("H61").Value = 2134.7222
$ws.Range("I61").Value = 2147.5715
$ws.Range("K61").Value = 2147.5715
$ws.Range("M61").Value = -1945.5715
$ws.Range("H93").Value = 3478.45
$ws.Range("I93").Value = 4300.4165
$ws.Range("J93").Value = 2245.5
$ws.Range("K93").Value = 4300.4165
$ws.Range("L93").Value = 2245.5
$ws.Range("M93").Value = -3052.4165
$ws.Range("N93").Value = -4741.5
$ws.Range("H113").Value = 2134.7222
$ws.Range("I113").Value = 2147.5715
$ws.Range("K113").Value = 2147.5715
$ws.Range("M113").Value = 22.42849999999999
$ws.Range("H122").Value = 3907.5
$ws.Range("I122").Value = 3436.818
$ws.Range("J122").Value = 5633.3335
$ws.Range("K122").Value = 10310.454
$ws.Range("L122").Value = 16900.0005
$ws.Range("M122").Value = -7860.454000000002
$ws.Range("N122").Value = -21800.0005
$ws.Range("H132").Value = 4294.793
$ws.Range("I132").Value = 4416.9546
$ws.Range("K132").Value = 13250.8638
$ws.Range("M132").Value = -10720.8638

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 942.7714
$ws.Range("I100").Value = 947.0714
$ws.Range("K100").Value = 1894.1428
$ws.Range("M100").Value = -1353.1428
$ws.Range("H107").Value = 762.35297
$ws.Range("I107").Value = 498.48
$ws.Range("J107").Value = 1495.3334
$ws.Range("K107").Value = 1495.44
$ws.Range("L107").Value = 4486.0002
$ws.Range("M107").Value = 424.5599999999999
$ws.Range("N107").Value = -8326.0002

